# Update rest-get.xlsx:
#  - Point the "findByTags" URL (used as the display text for cell E3's
#    hyperlink) at localhost instead of the old demo server.
#  - Reset the sheet view so it is no longer scrolled to column D and the
#    active selection sits on E3 instead of E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rest-Get")

# Cell E3 holds the URL text (and carries the hyperlink to the same address).
$ws.Range("E3").Value = "https://localhost/pets/findByTags?tags=red"

# Scroll the view back to the top-left (removes the topLeftCell="D1" pin)
# and move/select the active cell to E3.
$window = $excel.ActiveWindow
$window.ScrollColumn = 1
$window.ScrollRow = 1
$ws.Range("E3").Select()
